# Updates cryptos list values to match the latest scrape (GitHub Actions run).
# Numeric-looking "Price" strings are written with a leading apostrophe so Excel
# keeps them as text (matching the original workbook formatting) instead of
# auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '29.837.32'; ForceText = $false }
    @{ Cell = "D3"; Value = '1.869.30'; ForceText = $false }
    @{ Cell = "E3"; Value = '  -1.40%  '; ForceText = $false }
    @{ Cell = "D4"; Value = '0.9989'; ForceText = $true }
    @{ Cell = "E4"; Value = '  -0.17%  '; ForceText = $false }
    @{ Cell = "D5"; Value = '0.7420'; ForceText = $true }
    @{ Cell = "E5"; Value = '  -4.17%  '; ForceText = $false }
    @{ Cell = "D7"; Value = '0.9987'; ForceText = $true }
    @{ Cell = "D8"; Value = '0.3153'; ForceText = $true }
    @{ Cell = "E8"; Value = '  +0.48%  '; ForceText = $false }
    @{ Cell = "B9"; Value = 'Dogecoin'; ForceText = $false }
    @{ Cell = "C9"; Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; ForceText = $false }
    @{ Cell = "D9"; Value = '0.07122'; ForceText = $true }
    @{ Cell = "E9"; Value = '  -2.19%  '; ForceText = $false }
    @{ Cell = "B10"; Value = 'Solana'; ForceText = $false }
    @{ Cell = "C10"; Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; ForceText = $false }
    @{ Cell = "D10"; Value = '24.65'; ForceText = $true }
    @{ Cell = "E10"; Value = '  -4.36%  '; ForceText = $false }
    @{ Cell = "D11"; Value = '0.08391'; ForceText = $true }
    @{ Cell = "E11"; Value = '  -6.24%  '; ForceText = $false }
    @{ Cell = "D12"; Value = '0.7521'; ForceText = $true }
    @{ Cell = "D13"; Value = '5.423'; ForceText = $true }
    @{ Cell = "E13"; Value = '  -0.58%  '; ForceText = $false }
    @{ Cell = "D14"; Value = '1.859.75'; ForceText = $false }
    @{ Cell = "E14"; Value = '  +0.28%  '; ForceText = $false }
    @{ Cell = "D15"; Value = '92.61'; ForceText = $true }
    @{ Cell = "E15"; Value = '  -2.16%  '; ForceText = $false }
    @{ Cell = "D16"; Value = '29.836.13'; ForceText = $false }
    @{ Cell = "E16"; Value = '  -0.01%  '; ForceText = $false }
    @{ Cell = "D17"; Value = '6.023'; ForceText = $true }
    @{ Cell = "E17"; Value = '  -2.80%  '; ForceText = $false }
    @{ Cell = "E18"; Value = '  -3.09%  '; ForceText = $false }
    @{ Cell = "D19"; Value = '242.69'; ForceText = $true }
    @{ Cell = "E19"; Value = '  -1.48%  '; ForceText = $false }
    @{ Cell = "D20"; Value = '0.000007813'; ForceText = $true }
    @{ Cell = "E20"; Value = '  -0.97%  '; ForceText = $false }
    @{ Cell = "D21"; Value = '0.9985'; ForceText = $true }
    @{ Cell = "E21"; Value = '  -0.13%  '; ForceText = $false }
    @{ Cell = "D22"; Value = '2.114.96'; ForceText = $false }
    @{ Cell = "E22"; Value = '  +2.34%  '; ForceText = $false }
    @{ Cell = "D23"; Value = '7.931'; ForceText = $true }
    @{ Cell = "E23"; Value = '  -2.54%  '; ForceText = $false }
    @{ Cell = "D24"; Value = '0.9960'; ForceText = $true }
    @{ Cell = "E24"; Value = '  -0.48%  '; ForceText = $false }
    @{ Cell = "D25"; Value = '0.1569'; ForceText = $true }
    @{ Cell = "E25"; Value = '  -1.32%  '; ForceText = $false }
    @{ Cell = "D26"; Value = '9.290'; ForceText = $true }
    @{ Cell = "E26"; Value = '  -2.70%  '; ForceText = $false }
    @{ Cell = "D27"; Value = '163.72'; ForceText = $true }
    @{ Cell = "E27"; Value = '  +0.29%  '; ForceText = $false }
    @{ Cell = "E28"; Value = '  -1.49%  '; ForceText = $false }
    @{ Cell = "D29"; Value = '2.020'; ForceText = $true }
    @{ Cell = "E29"; Value = '  -1.40%  '; ForceText = $false }
    @{ Cell = "D30"; Value = '1.470'; ForceText = $true }
    @{ Cell = "E30"; Value = '  +3.03%  '; ForceText = $false }
    @{ Cell = "D31"; Value = '4.636'; ForceText = $true }
    @{ Cell = "E31"; Value = '  +1.91%  '; ForceText = $false }
    @{ Cell = "D32"; Value = '1.530'; ForceText = $true }
    @{ Cell = "E32"; Value = '  -0.98%  '; ForceText = $false }
    @{ Cell = "D33"; Value = '4.294'; ForceText = $true }
    @{ Cell = "E33"; Value = '  +4.25%  '; ForceText = $false }
    @{ Cell = "D34"; Value = '0.05309'; ForceText = $true }
    @{ Cell = "E34"; Value = '  -3.93%  '; ForceText = $false }
    @{ Cell = "E35"; Value = '  -1.37%  '; ForceText = $false }
    @{ Cell = "D36"; Value = '0.7511'; ForceText = $true }
    @{ Cell = "E36"; Value = '  -0.35%  '; ForceText = $false }
    @{ Cell = "D37"; Value = '0.9997'; ForceText = $true }
    @{ Cell = "E37"; Value = '  +0.32%  '; ForceText = $false }
    @{ Cell = "E38"; Value = '  -0.85%  '; ForceText = $false }
    @{ Cell = "D39"; Value = '0.01953'; ForceText = $true }
    @{ Cell = "E39"; Value = '  -0.51%  '; ForceText = $false }
    @{ Cell = "D40"; Value = '2.743'; ForceText = $true }
    @{ Cell = "E40"; Value = '  -1.78%  '; ForceText = $false }
    @{ Cell = "D41"; Value = '0.4466'; ForceText = $true }
    @{ Cell = "E41"; Value = '  -1.09%  '; ForceText = $false }
    @{ Cell = "D42"; Value = '1.105.42'; ForceText = $false }
    @{ Cell = "E42"; Value = '  +1.64%  '; ForceText = $false }
    @{ Cell = "D43"; Value = '6.079'; ForceText = $true }
    @{ Cell = "E43"; Value = '  +0.37%  '; ForceText = $false }
    @{ Cell = "D44"; Value = '72.04'; ForceText = $true }
    @{ Cell = "E44"; Value = '  -2.68%  '; ForceText = $false }
    @{ Cell = "D45"; Value = '0.8600'; ForceText = $true }
    @{ Cell = "E45"; Value = '  +0.54%  '; ForceText = $false }
    @{ Cell = "D46"; Value = '0.9994'; ForceText = $true }
    @{ Cell = "E46"; Value = '  -0.08%  '; ForceText = $false }
    @{ Cell = "D47"; Value = '103.21'; ForceText = $true }
    @{ Cell = "E47"; Value = '  +0.45%  '; ForceText = $false }
    @{ Cell = "D48"; Value = '7.674'; ForceText = $true }
    @{ Cell = "E48"; Value = '  +0.74%  '; ForceText = $false }
    @{ Cell = "D49"; Value = '1.839'; ForceText = $true }
    @{ Cell = "E49"; Value = '  -3.06%  '; ForceText = $false }
    @{ Cell = "D50"; Value = '3.036'; ForceText = $true }
    @{ Cell = "E50"; Value = '  +1.21%  '; ForceText = $false }
    @{ Cell = "D51"; Value = '2.015.08'; ForceText = $false }
    @{ Cell = "E51"; Value = '  -0.47%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $v = $u.Value
    if ($u.ForceText) {
        # Leading apostrophe tells Excel to store this as text, not a number
        $v = "'" + $v
    }
    $ws.Range($u.Cell).Value = $v
}
